# Bill of Materials - adapt resistance for production
#
# The BOM row that used to lump the 68-ohm resistor group together with
# two mis-classified 180-ohm resistors (R107, R137) is split into two
# separate rows:
#   - row 39 keeps the 12 genuine 68-ohm resistors (value corrected
#     from "180"/mixed to "68")
#   - a new row 40 is inserted holding the two 180-ohm resistors
#     (R107, R137) that were incorrectly folded into the 68-ohm row
#
# All of the following rows shift down by one. The script also cleans up
# left-over "left align" styling on column B (and the stale selection
# formatting on the old row 47), widens column C, and resets the
# worksheet view back to the top-left corner.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row at position 40, pushing rows 40-67 down to 41-68.
$ws.Rows.Item(40).Insert()

# 2) Correct row 39: the 68-ohm resistor group, minus R107/R137.
$ws.Cells.Item(39, 1).Value = 12
$ws.Cells.Item(39, 2).Value = 68
$ws.Cells.Item(39, 3).Value = "R68, R69, R70, R71, R72, R73, R98, R99, R100, R147, R148, R149"
$ws.Cells.Item(39, 4).Value = 68
$ws.Cells.Item(39, 5).Value = "RESC1608X06L"

# 3) New row 40: the 180-ohm resistor group (R107, R137).
$ws.Cells.Item(40, 1).Value = 2
$ws.Cells.Item(40, 2).Value = 180
$ws.Cells.Item(40, 3).Value = "R107, R137"
$ws.Cells.Item(40, 4).Value = "Resistor"
$ws.Cells.Item(40, 5).Value = "RESC1608X06L"

# 4) Remove the left-alignment style from column B (now general, s=0),
#    and from the formerly-special A47:C47 (now A48:C48 after the
#    insert) selection formatting.
$ws.Range("B1:B68").HorizontalAlignment = 1
$ws.Range("A48:C48").HorizontalAlignment = 1

# 5) Column widths: C grows, B stays the same.
$ws.Columns.Item(3).ColumnWidth = 204.15

# 6) Reset the view to the top-left corner with a simple A1 selection.
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("A1").Select()
